$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).ClearContents()

$ws.Range("C1").Select()
